# Add Q3-2022 data: a new "2022-Q3" worksheet (copied from "2022-Q2" so it
# keeps the same layout/styles) placed right before "2022-Q2", plus a new
# summary row on "总计" for the new quarter. All older quarters simply shift
# position (their own data is untouched).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (totals) sheet: push the existing data rows down one slot and
#    insert the new 2022-Q3 totals at the top of the data (row 2).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Duplicate the last row (old row 4, "2021-Q3") into the new row 5,
# copying the cell first so A5 keeps the same number style as A2:A4.
$wsTotal.Cells.Item(4, 1).Copy($wsTotal.Cells.Item(5, 1))
$wsTotal.Cells.Item(5, 1).Value2 = 3
$wsTotal.Cells.Item(5, 2).Value2 = "2021-Q3"
$wsTotal.Cells.Item(5, 3).Value2 = 2
$wsTotal.Cells.Item(5, 4).Value2 = 0.6

# Shift the quarter labels / values down one row, in place.
$wsTotal.Cells.Item(4, 2).Value2 = "2021-Q4"
$wsTotal.Cells.Item(4, 4).Value2 = 0.66

$wsTotal.Cells.Item(3, 2).Value2 = "2022-Q2"
$wsTotal.Cells.Item(3, 4).Value2 = 2.39

$wsTotal.Cells.Item(2, 2).Value2 = "2022-Q3"
$wsTotal.Cells.Item(2, 4).Value2 = 2.98

# ---------------------------------------------------------------------
# 2. New "2022-Q3" worksheet: duplicate "2022-Q2" (same headers/styles)
#    and place the copy right before it, then overwrite with the new
#    quarter's fund data.
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

$wsQ3.Cells.Item(2, 2).Value2 = "'012349"
$wsQ3.Cells.Item(2, 3).Value2 = "天弘恒生科技指数（QDII）C"
$wsQ3.Cells.Item(2, 4).Value2 = "'33.57"
$wsQ3.Cells.Item(2, 5).Value2 = "'92.84"
$wsQ3.Cells.Item(2, 6).Value2 = "'4.64"
$wsQ3.Cells.Item(2, 7).Value2 = "'1.5576"
$wsQ3.Cells.Item(2, 8).Value2 = 9

$wsQ3.Cells.Item(3, 2).Value2 = "'012348"
$wsQ3.Cells.Item(3, 3).Value2 = "天弘恒生科技指数（QDII）A"
$wsQ3.Cells.Item(3, 4).Value2 = "'30.64"
$wsQ3.Cells.Item(3, 5).Value2 = "'92.84"
$wsQ3.Cells.Item(3, 6).Value2 = "'4.64"
$wsQ3.Cells.Item(3, 7).Value2 = "'1.4217"
$wsQ3.Cells.Item(3, 8).Value2 = 9
